# Add two new columns "I0" (column I) and "IF" (column J) to the worksheet,
# matching the header style already used by the other header cells (e.g. H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy formatting from the existing header cell H1 so the
# new header cells reuse the same style (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-14
$I = @(8, 8, 6, 5, 7, 4, 8, 6, 4, 4, 7, 7, 5)
$J = @(8, 8, 6, 5, 8, 5, 8, 6, 4, 4, 7, 7, 5)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
